$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the sheet completely so stale cells (with stale styles) are not left behind
for ($i = 0; $i -lt 26; $i++) {
    $ws.Rows.Item(1).Delete()
}

$c = $ws.Range('B1')
$c.Value = 'Ementa atual:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C1')
$c.Value = 'Ementa modificada (dados modificados em vermelho):'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B2')
$c.Value = 'LOM3043'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C2')
$c.Value = 'LOM3043'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$c = $ws.Range('A3')
$c.Value = 'Nome:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B3')
$c.Value = ' Seleção de Materiais'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C3')
$c.Value = ' Seleção de Materiais'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$c = $ws.Range('A4')
$c.Value = 'Name:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B4')
$c.Value = 'Selection of Materials'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C4')
$c.Value = 'Selection of Materials'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$c = $ws.Range('A5')
$c.Value = 'Créditos-aula:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B5')
$c.Value = '2'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C5')
$c.Value = '2'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$c = $ws.Range('A6')
$c.Value = 'Créditos-trabalho'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B6')
$c.Value = '0'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C6')
$c.Value = '0'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$c = $ws.Range('A7')
$c.Value = 'Carga horária:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B7')
$c.Value = '30 h'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C7')
$c.Value = '30 h'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$c = $ws.Range('A8')
$c.Value = 'Ativação:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B8')
$c.Value = '01/01/2020'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C8')
$c.Value = '01/01/2020'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$c = $ws.Range('A9')
$c.Value = 'Semestre ideal:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B9')
$c.Value = 'EM-8'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C9')
$c.Value = 'EM-8'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$c = $ws.Range('A10')
$c.Value = 'Objetivos:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B10')
$c.Value = '5840622 - Miguel Justino Ribeiro Barboza'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C10')
$c.Value = '5840622 - Miguel Justino Ribeiro Barboza'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$ws.Rows.Item(10).RowHeight = 60

$c = $ws.Range('A11')
$c.Value = 'Objectives:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$ws.Rows.Item(11).RowHeight = 60

$c = $ws.Range('A12')
$c.Value = 'Docentes responsáveis:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('A13')
$c.Value = 'Programa resumido:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B13')
$c.Value = 'Semestral'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C13')
$c.Value = 'Semestral'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$ws.Rows.Item(13).RowHeight = 60

$c = $ws.Range('A14')
$c.Value = 'Short syllabus:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$ws.Rows.Item(14).RowHeight = 60

$c = $ws.Range('A15')
$c.Value = 'Programa:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B15')
$c.Value = '01/01/2020'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C15')
$c.Value = '01/01/2020'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$ws.Rows.Item(15).RowHeight = 120

$c = $ws.Range('A16')
$c.Value = 'Syllabus:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$ws.Rows.Item(16).RowHeight = 120

$c = $ws.Range('A17')
$c.Value = 'Avaliação:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('A18')
$c.Value = 'Método:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B18')
$c.Value = '5840622 - Miguel Justino Ribeiro Barboza'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C18')
$c.Value = '5840622 - Miguel Justino Ribeiro Barboza'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$ws.Rows.Item(18).RowHeight = 60

$c = $ws.Range('A19')
$c.Value = 'Critério:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B19')
$c.Value = 'Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa.'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C19')
$c.Value = 'Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa.'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$ws.Rows.Item(19).RowHeight = 60

$c = $ws.Range('A20')
$c.Value = 'Norma de recuperação:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B20')
$c.Value = 'A média do semestre será computada com base na relação:M=(P1+2P2)/3'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C20')
$c.Value = 'A média do semestre será computada com base na relação:M=(P1+2P2)/3'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$ws.Rows.Item(20).RowHeight = 60

$c = $ws.Range('A21')
$c.Value = 'Bibliografia:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B21')
$c.Value = 'A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C21')
$c.Value = 'A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$ws.Rows.Item(21).RowHeight = 120

$c = $ws.Range('A22')
$c.Value = 'Requisitos:'
$c.Font.Bold = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('B23')
$c.Value = 'LOM3036 -  Propriedades Mecânicas  (Requisito fraco)
'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C23')
$c.Value = 'LOM3036 -  Propriedades Mecânicas  (Requisito fraco)
'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$ws.Rows.Item(23).RowHeight = 30

$c = $ws.Range('B24')
$c.Value = 'LOM3057 -  Introdução aos Materiais Poliméricos  (Requisito fraco)
'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C24')
$c.Value = 'LOM3057 -  Introdução aos Materiais Poliméricos  (Requisito fraco)
'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$ws.Rows.Item(24).RowHeight = 30

$c = $ws.Range('B25')
$c.Value = 'LOM3082 -  Cerâmica Física  (Requisito fraco)
'
$c.WrapText = $true
$c.VerticalAlignment = -4160

$c = $ws.Range('C25')
$c.Value = 'LOM3082 -  Cerâmica Física  (Requisito fraco)
'
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c.Font.Color = 255

$ws.Rows.Item(25).RowHeight = 30
